$d = $word.ActiveDocument

# Locate the "August, 2017" paragraph (the thesis title-page date line) and the
# trailing paragraphs that follow it: the "Statement of Originality" heading,
# the originality-statement body text, and the "Signed: ____" line (plus the
# blank spacer paragraphs interleaved between them).
$augustPara = $null
$stmtPara = $null
$signedPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "August, 2017") {
        $augustPara = $p
    }
    if ($t -match "Statement of Originality") {
        $stmtPara = $p
    }
    if ($t -match "Signed:") {
        $signedPara = $p
    }
}

if ($augustPara -eq $null) { throw "could not find the 'August, 2017' paragraph" }
if ($stmtPara -eq $null) { throw "could not find the 'Statement of Originality' paragraph" }
if ($signedPara -eq $null) { throw "could not find the 'Signed:' paragraph" }

# Remove the whole "Statement of Originality" section: the heading paragraph
# through the final "Signed: ____" paragraph, inclusive (this deletes the
# page break that introduced the heading too, since it lives inside that
# paragraph's range).
$deleteRange = $d.Range($stmtPara.Range.Start, $signedPara.Range.End)
$deleteRange.Delete()

# Rebuild the "August, 2017" paragraph in place: bump the paragraph-mark's
# run size to 48, split the run text into "August," / " 2017" wrapped around
# a pair of proofErr grammar-check markers, and re-home the _GoBack bookmark
# at the end of the paragraph (it used to sit inside the now-deleted text).
$augustXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00FC08B9" w:rsidRPr="00685DF6" w:rsidRDefault="00FC08B9" w:rsidP="00685DF6"><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="48"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="36"/></w:rPr><w:t>August,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve"> 2017</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$augustPara.Range.InsertXML($augustXml)
